{"js": "// Remove the \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer line, and the\n// blank paragraph that directly followed the footer line -- this is the\n// static-site boilerplate block that was dropped from the page at build\n// time. The blank paragraph that precedes \"Ver no Jupiter...\" (right after\n// the bibliography text) is left untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightMarker = \"Contact: luizeleno@usp.br\";\n\nconst toDelete = [];\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === jupiterText) {\n    toDelete.push(i);\n  } else if (t.indexOf(copyrightMarker) !== -1) {\n    toDelete.push(i);\n    // Blank spacer paragraph right after the footer line goes too.\n    if (i + 1 < items.length && items[i + 1].text.trim().length === 0) {\n      toDelete.push(i + 1);\n    }\n  }\n}\n\n// Delete from the end first so earlier indices stay valid.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer line, and the\n# blank paragraph that directly followed the footer line -- this is the\n# static-site boilerplate block that was dropped from the page at build\n# time. The blank paragraph that precedes \"Ver no Jupiter...\" (right after\n# the bibliography text) is left untouched.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightPrefix = \"Contact: luizeleno@usp.br\"\n\n# Walk backwards so deleting a paragraph never shifts the index of one we\n# still need to visit.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n\n    if ($t.Contains($copyrightPrefix)) {\n        # Blank spacer paragraph right after the footer line goes too.\n        $next = $d.Paragraphs.Item($i + 1)\n        if ($next.Range.Text.Trim().Length -eq 0) {\n            $next.Range.Delete()\n        }\n        $p.Range.Delete()\n    }\n    elseif ($t -eq $jupiterText) {\n        $p.Range.Delete()\n    }\n}\n\nWrite-Output \"done\"\n"}
